# Updates cryptos list data (Price and Volume(1h) columns) on worksheet "sheet1"
# to match the latest scrape, per the commit:
#   "Updated cryptos list on Sat May 27 07:16:24 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
# Values that look like plain numbers (e.g. "309.28") are written with a
# leading apostrophe so Excel keeps them as literal text (matching the
# workbook author's original inline-string cells) instead of turning them
# into numeric cells.
$updates = @(
    @{ Addr = "D2"; Value = '26.922.04' },
    @{ Addr = "E2"; Value = '  +1.12%  ' },
    @{ Addr = "D3"; Value = '1.845.62' },
    @{ Addr = "E3"; Value = '  +1.17%  ' },
    @{ Addr = "E4"; Value = '  +0.03%  ' },
    @{ Addr = "D5"; Value = '''309.28' },
    @{ Addr = "E5"; Value = '  +0.73%  ' },
    @{ Addr = "E6"; Value = '  -0.04%  ' },
    @{ Addr = "D7"; Value = '''0.4760' },
    @{ Addr = "E7"; Value = '  +2.47%  ' },
    @{ Addr = "D8"; Value = '''0.3667' },
    @{ Addr = "E8"; Value = '  +1.76%  ' },
    @{ Addr = "D9"; Value = '''0.07198' },
    @{ Addr = "E9"; Value = '  +0.87%  ' },
    @{ Addr = "D10"; Value = '''0.9266' },
    @{ Addr = "E10"; Value = '  +2.79%  ' },
    @{ Addr = "D11"; Value = '''19.73' },
    @{ Addr = "E11"; Value = '  +1.70%  ' },
    @{ Addr = "D12"; Value = '''0.07683' },
    @{ Addr = "E12"; Value = '  -1.09%  ' },
    @{ Addr = "D13"; Value = '1.826.68' },
    @{ Addr = "E13"; Value = '  +1.81%  ' },
    @{ Addr = "D14"; Value = '''5.311' },
    @{ Addr = "E14"; Value = '  +0.95%  ' },
    @{ Addr = "D15"; Value = '''6.408' },
    @{ Addr = "E15"; Value = '  +1.24%  ' },
    @{ Addr = "D16"; Value = '''88.69' },
    @{ Addr = "E16"; Value = '  +1.50%  ' },
    @{ Addr = "D17"; Value = '''1.009' },
    @{ Addr = "E17"; Value = '  +0.02%  ' },
    @{ Addr = "D18"; Value = '''0.000008625' },
    @{ Addr = "E18"; Value = '  +0.72%  ' },
    @{ Addr = "E19"; Value = '  -0.01%  ' },
    @{ Addr = "D20"; Value = '26.954.61' },
    @{ Addr = "E20"; Value = '  +1.11%  ' },
    @{ Addr = "D21"; Value = '''14.54' },
    @{ Addr = "E21"; Value = '  +2.76%  ' },
    @{ Addr = "E22"; Value = '  +0.72%  ' },
    @{ Addr = "D23"; Value = '''10.62' },
    @{ Addr = "E23"; Value = '  +0.94%  ' },
    @{ Addr = "D24"; Value = '''1.921' },
    @{ Addr = "E24"; Value = '  +0.09%  ' },
    @{ Addr = "D25"; Value = '''152.36' },
    @{ Addr = "E25"; Value = '  +0.04%  ' },
    @{ Addr = "D26"; Value = '''18.14' },
    @{ Addr = "D27"; Value = '''1.997' },
    @{ Addr = "E27"; Value = '  +1.10%  ' },
    @{ Addr = "D28"; Value = '''114.23' },
    @{ Addr = "E28"; Value = '  +0.45%  ' },
    @{ Addr = "D29"; Value = '''4.920' },
    @{ Addr = "E29"; Value = '  +2.39%  ' },
    @{ Addr = "D30"; Value = '''0.08883' },
    @{ Addr = "E30"; Value = '  +0.87%  ' },
    @{ Addr = "D31"; Value = '''3.315' },
    @{ Addr = "E31"; Value = '  +5.59%  ' },
    @{ Addr = "D32"; Value = '''1.174' },
    @{ Addr = "E32"; Value = '  +2.97%  ' },
    @{ Addr = "D33"; Value = '''0.7458' },
    @{ Addr = "E33"; Value = '  +1.97%  ' },
    @{ Addr = "D34"; Value = '''4.481' },
    @{ Addr = "E34"; Value = '  +0.94%  ' },
    @{ Addr = "D35"; Value = '''2.737' },
    @{ Addr = "E35"; Value = '  +0.58%  ' },
    @{ Addr = "D36"; Value = '''1.102' },
    @{ Addr = "E36"; Value = '  +2.58%  ' },
    @{ Addr = "D37"; Value = '''0.01955' },
    @{ Addr = "E37"; Value = '  +1.49%  ' },
    @{ Addr = "D38"; Value = '''0.05264' },
    @{ Addr = "E38"; Value = '  +2.94%  ' },
    @{ Addr = "D39"; Value = '''2.972' },
    @{ Addr = "E39"; Value = '  +1.59%  ' },
    @{ Addr = "D40"; Value = '''0.5192' },
    @{ Addr = "E40"; Value = '  +2.71%  ' },
    @{ Addr = "D41"; Value = '''6.960' },
    @{ Addr = "E41"; Value = '  +0.76%  ' },
    @{ Addr = "D42"; Value = '''0.1512' },
    @{ Addr = "E42"; Value = '  +1.13%  ' },
    @{ Addr = "D43"; Value = '''8.210' },
    @{ Addr = "E43"; Value = '  +2.68%  ' },
    @{ Addr = "D44"; Value = '''10.54' },
    @{ Addr = "E44"; Value = '  +5.44%  ' },
    @{ Addr = "D45"; Value = '''0.4715' },
    @{ Addr = "E45"; Value = '  +1.30%  ' },
    @{ Addr = "E46"; Value = '  +0.02%  ' },
    @{ Addr = "D47"; Value = '''101.50' },
    @{ Addr = "E47"; Value = '  +3.65%  ' },
    @{ Addr = "E48"; Value = '  +3.01%  ' },
    @{ Addr = "D49"; Value = '''66.13' },
    @{ Addr = "E49"; Value = '  +3.88%  ' },
    @{ Addr = "D50"; Value = '''0.06021' },
    @{ Addr = "E50"; Value = '  +0.65%  ' },
    @{ Addr = "D51"; Value = '''0.8847' },
    @{ Addr = "E51"; Value = '  +3.89%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# Cells that were forced to text via a leading apostrophe pick up an
# implicit "quote prefix" style; reset them back to the default "Normal"
# style so formatting is unaffected, matching the original (unstyled) cells.
foreach ($u in $updates) {
    if ($u.Value.StartsWith("'")) {
        $ws.Range($u.Addr).Style = "Normal"
    }
}
